$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set status "Not Started" and member "Jun" for the newly-assigned EDA rows
$statusCells = @("E8", "E33", "E34", "E35", "E36", "E37", "E53")
foreach ($cell in $statusCells) {
    $ws.Range($cell).Value = "Not Started"
}

$memberCells = @("F33", "F34", "F35", "F36", "F37", "F53")
foreach ($cell in $memberCells) {
    $ws.Range($cell).Value = "Jun"
}

# Reflect the updated selection state recorded in the saved file
$ws.Range("E37").Select()
